$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H55").Value = 90910210
$ws.Range("I55").Value = 142858700
$ws.Range("J55").Value = 345
$ws.Range("K55").Value = 142858700
$ws.Range("L55").Value = 345
$ws.Range("M55").Value = -142858486
$ws.Range("N55").Value = -773
$ws.Range("H137").Value = 1426.6842
$ws.Range("I137").Value = 1567.1111
$ws.Range("J137").Value = 1300.3
$ws.Range("K137").Value = 4701.3333
$ws.Range("L137").Value = 3900.9
$ws.Range("M137").Value = -2151.3333
$ws.Range("N137").Value = -9000.9
$ws.Range("H138").Value = 2709.0334
$ws.Range("I138").Value = 2577
$ws.Range("J138").Value = 2753.0444
$ws.Range("K138").Value = 7731
$ws.Range("L138").Value = 8259.1332
$ws.Range("M138").Value = -2591
$ws.Range("N138").Value = -18539.1332
$ws.Range("H141").Value = 4319.033
$ws.Range("I141").Value = 1789.9048
$ws.Range("K141").Value = 5369.7144
$ws.Range("M141").Value = -189.7143999999998

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 3325.6365
$ws.Range("I45").Value = 2712.4
$ws.Range("J45").Value = 3836.6667
$ws.Range("K45").Value = 2712.4
$ws.Range("L45").Value = 3836.6667
$ws.Range("M45").Value = -2335.4
$ws.Range("N45").Value = -4590.6667
$ws.Range("H55").Value = 29599.2
$ws.Range("J55").Value = 29599.2
$ws.Range("L55").Value = 29599.2
$ws.Range("N55").Value = -30229.2
$ws.Range("H61").Value = 8335370.5
$ws.Range("I61").Value = 15152822
$ws.Range("J61").Value = 2929.3333
$ws.Range("K61").Value = 15152822
$ws.Range("L61").Value = 2929.3333
$ws.Range("M61").Value = -15152610
$ws.Range("N61").Value = -3353.3333
$ws.Range("H74").Value = 725.3611
$ws.Range("I74").Value = 445.75
$ws.Range("J74").Value = 1074.875
$ws.Range("K74").Value = 445.75
$ws.Range("L74").Value = 1074.875
$ws.Range("M74").Value = 428.25
$ws.Range("N74").Value = -2822.875
$ws.Range("H77").Value = 725.3611
$ws.Range("I77").Value = 445.75
$ws.Range("J77").Value = 1074.875
$ws.Range("K77").Value = 2228.75
$ws.Range("L77").Value = 5374.375
$ws.Range("M77").Value = 2139.25
$ws.Range("N77").Value = -14110.375
$ws.Range("H136").Value = 8335370.5
$ws.Range("I136").Value = 15152822
$ws.Range("J136").Value = 2929.3333
$ws.Range("K136").Value = 45458466
$ws.Range("L136").Value = 8787.999899999999
$ws.Range("M136").Value = -45455916
$ws.Range("N136").Value = -13887.9999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 2770.3428
$ws.Range("I134").Value = 2574
$ws.Range("J134").Value = 3146.6667
$ws.Range("K134").Value = 7722
$ws.Range("L134").Value = 9440.000100000001
$ws.Range("M134").Value = -5187
$ws.Range("N134").Value = -14510.0001
$ws.Range("H135").Value = 49312.5
$ws.Range("I135").Value = 0
$ws.Range("J135").Value = 49312.5
$ws.Range("K135").Value = 0
$ws.Range("L135").Value = 49312.5
$ws.Range("M135").ClearContents()
$ws.Range("N135").Value = -59452.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4514.15
$ws.Range("I31").Value = 1563.3334
$ws.Range("J31").Value = 6284.64
$ws.Range("K31").Value = 1563.3334
$ws.Range("L31").Value = 6284.64
$ws.Range("M31").Value = -1268.3334
$ws.Range("N31").Value = -6874.64
$ws.Range("H34").Value = 4514.15
$ws.Range("I34").Value = 1563.3334
$ws.Range("J34").Value = 6284.64
$ws.Range("K34").Value = 1563.3334
$ws.Range("L34").Value = 6284.64
$ws.Range("M34").Value = -1361.3334
$ws.Range("N34").Value = -6688.64
$ws.Range("H45").Value = 0
$ws.Range("J45").Value = 0
$ws.Range("L45").Value = 0
$ws.Range("N45").ClearContents()
$ws.Range("H58").Value = 2909.8125
$ws.Range("I58").Value = 2812.1538
$ws.Range("J58").Value = 3333
$ws.Range("K58").Value = 2812.1538
$ws.Range("L58").Value = 3333
$ws.Range("M58").Value = -2609.1538
$ws.Range("N58").Value = -3739
$ws.Range("H132").Value = 6668902.5
$ws.Range("I132").Value = 1744.9231
$ws.Range("K132").Value = 5234.7693
$ws.Range("M132").Value = -2704.7693
$ws.Range("H134").Value = 1512.0834
$ws.Range("I134").Value = 855.1667
$ws.Range("J134").Value = 2169
$ws.Range("K134").Value = 2565.5001
$ws.Range("L134").Value = 6507
$ws.Range("M134").Value = -30.5001000000002
$ws.Range("N134").Value = -11577
$ws.Range("H136").Value = 2909.8125
$ws.Range("I136").Value = 2812.1538
$ws.Range("J136").Value = 3333
$ws.Range("K136").Value = 8436.4614
$ws.Range("L136").Value = 9999
$ws.Range("M136").Value = -5886.4614
$ws.Range("N136").Value = -15099

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 12402.167
$ws.Range("I3").Value = 8625
$ws.Range("J3").Value = 19956.5
$ws.Range("K3").Value = 25875
$ws.Range("L3").Value = 59869.5
$ws.Range("M3").Value = -25763
$ws.Range("N3").Value = -60093.5
$ws.Range("H68").Value = 1349.6666
$ws.Range("I68").Value = 774.0714
$ws.Range("J68").Value = 1586.6765
$ws.Range("K68").Value = 2322.2142
$ws.Range("L68").Value = 4760.029500000001
$ws.Range("M68").Value = -1511.2142
$ws.Range("N68").Value = -6382.029500000001
$ws.Range("H71").Value = 1349.6666
$ws.Range("I71").Value = 774.0714
$ws.Range("J71").Value = 1586.6765
$ws.Range("K71").Value = 6966.6426
$ws.Range("L71").Value = 14280.0885
$ws.Range("M71").Value = -2910.6426
$ws.Range("N71").Value = -22392.0885
$ws.Range("H114").Value = 921
$ws.Range("I114").Value = 260
$ws.Range("J114").Value = 1317.6
$ws.Range("K114").Value = 780
$ws.Range("L114").Value = 3952.8
$ws.Range("M114").Value = 2474
$ws.Range("N114").Value = -10460.8
$ws.Range("H131").Value = 1142.091
$ws.Range("I131").Value = 943.3333
$ws.Range("J131").Value = 1173.4736
$ws.Range("K131").Value = 2829.9999
$ws.Range("L131").Value = 3520.4208
$ws.Range("M131").Value = 2210.0001
$ws.Range("N131").Value = -13600.4208

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 90912260
$ws.Range("I7").Value = 166669300
$ws.Range("J7").Value = 3801
$ws.Range("K7").Value = 166669300
$ws.Range("L7").Value = 3801
$ws.Range("M7").Value = -166669188
$ws.Range("N7").Value = -4025
$ws.Range("H40").Value = 111114000
$ws.Range("I40").Value = 200002300
$ws.Range("J40").Value = 3625
$ws.Range("K40").Value = 200002300
$ws.Range("L40").Value = 3625
$ws.Range("M40").Value = -200002164
$ws.Range("N40").Value = -3897
$ws.Range("H126").Value = 90912260
$ws.Range("I126").Value = 166669300
$ws.Range("J126").Value = 3801
$ws.Range("K126").Value = 500007900
$ws.Range("L126").Value = 11403
$ws.Range("M126").Value = -500005430
$ws.Range("N126").Value = -16343
$ws.Range("H132").Value = 2840.3333
$ws.Range("I132").Value = 2733.9583
$ws.Range("J132").Value = 2961.9048
$ws.Range("K132").Value = 8201.874899999999
$ws.Range("L132").Value = 8885.714399999999
$ws.Range("M132").Value = -5671.874899999999
$ws.Range("N132").Value = -13945.7144
$ws.Range("H136").Value = 5209457
$ws.Range("I136").Value = 1068
$ws.Range("K136").Value = 3204
$ws.Range("M136").Value = -654

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H47").Value = 50000
$ws.Range("J47").Value = 50000
$ws.Range("L47").Value = 50000
$ws.Range("N47").Value = -51144
$ws.Range("H132").Value = 4558870.5
$ws.Range("I132").Value = 1400.1143
$ws.Range("J132").Value = 10059265
$ws.Range("K132").Value = 4200.3429
$ws.Range("L132").Value = 30177795
$ws.Range("M132").Value = -1670.3429
$ws.Range("N132").Value = -30182855
$ws.Range("H136").Value = 2359.6667
$ws.Range("I136").Value = 2226.9756
$ws.Range("J136").Value = 2699.6875
$ws.Range("K136").Value = 6680.926800000001
$ws.Range("L136").Value = 8099.0625
$ws.Range("M136").Value = -4130.926800000001
$ws.Range("N136").Value = -13199.0625
